# [REF] Tools refactoring step 9
#
# The "account_tax" sheet is missing values in column E ("amount") for a
# number of tax rows. Backfill those gaps with 0, matching the numeric
# values already present on the surrounding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("account_tax")
$ws.Activate()

# Rows whose column E cell is currently empty and needs to be populated
# with a numeric 0, same as the other (already filled-in) rows.
$targetRange = $ws.Range("E8,E9,E22:E43,E46:E57")
foreach ($area in $targetRange.Areas) {
    $area.Value = 0
}

# Leave the sheet's selection/cursor on E1, matching the final state of
# the workbook after the edit.
$ws.Range("A1").Select()
$ws.Range("E1").Select()
